$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume 1h (E) values for crypto rows
$ws.Range("D2").Value = "42.554.87"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "2.512.74"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "314.59"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("D6").Value = "94.22"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("D7").Value = "0.578"
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "35.63"
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "7.51"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D14").Value = "2.897.13"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "42.646.52"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").Value = "12.92"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("D21").Value = "0.0₃0958"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").Value = "69.29"
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("D23").Value = "250.23"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").Value = "2.01"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "26.65"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D29").Value = "40.84"
$ws.Range("E29").Value = "  +7.68%  "
$ws.Range("D30").Value = "10.26"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "5.92"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").Value = "156.14"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "19.19"
$ws.Range("E33").Value = "  +3.89%  "
$ws.Range("D35").Value = "3.27"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "0.0781"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D40").Value = "23.70"
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("D41").Value = "2.29"
$ws.Range("E41").Value = "  +13.35%  "
$ws.Range("D43").Value = "0.0303"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").Value = "3.76"
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("D45").Value = "3.30"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "2.013.71"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").Value = "85.42"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").Value = "8.80"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "2.751.49"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").Value = "73.15"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("D51").Value = "102.28"
$ws.Range("E51").Value = "  +0.99%  "

# Update Volume 1h (E) only for rows where Price stayed the same
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E13").Value = "  -3.91%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  +4.44%  "
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("E38").Value = "  -3.78%  "
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("E42").Value = "  +0.06%  "

# Rows 15 and 16 swap: Chainlink and WrappedEther exchange positions with updated values
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.516.48"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "15.11"
$ws.Range("E16").Value = "  -0.30%  "
